$p = $ppt.ActivePresentation

# "Title and Content" is CustomLayout #2 - it has a title placeholder and a
# body/content placeholder (idx=1), matching the placeholders used by the
# two new slides.
$layout = $p.SlideMaster.CustomLayouts.Item(2)

# ======================================================================
# New slide 2: outline slide ("Define GIT" / workflow / solo work / ...)
# ======================================================================
$s2 = $p.Slides.AddSlide(2, $layout)

$tr = $s2.Shapes.Item(2).TextFrame.TextRange

$tr.Text = "Define GIT"
$tr.InsertAfter("`rWork flow (analogy to keeping a note book)") | Out-Null
$tr.InsertAfter("`rDoing solo work ") | Out-Null

$r4a = $tr.InsertAfter("`rMaintaining your local repository (add, ")
$r4b = $r4a.InsertAfter("gitignore")
$r4b.InsertAfter(", revert, branch, checkout, remove, copy, merge, log, commit)") | Out-Null

$r5a = $tr.InsertAfter("`rCollaboration: Pushing to a remote repository in ")
$r5b = $r5a.InsertAfter("Github")
$r5b.InsertAfter("/cluster (push, merge, branch)") | Out-Null

$r6a = $tr.InsertAfter("`rDemo ")
$r6a.InsertAfter("for the code") | Out-Null

$tr.InsertAfter("`r") | Out-Null
$tr.InsertAfter("`r") | Out-Null
$tr.InsertAfter("`r") | Out-Null
$tr.InsertAfter("`r") | Out-Null

# Sub-bullets (level 2 in the UI == IndentLevel 2 == lvl="1" in the XML).
$tr.Paragraphs(4).IndentLevel = 2
$tr.Paragraphs(5).IndentLevel = 2
$tr.Paragraphs(7).IndentLevel = 2
$tr.Paragraphs(8).IndentLevel = 2
$tr.Paragraphs(9).IndentLevel = 2

# Final (empty) paragraph: plain text, no bullet.
$tr.Paragraphs(10).ParagraphFormat.Bullet.Visible = 0

# ======================================================================
# New slide 3: blank outline slide (title + content placeholder, empty)
# ======================================================================
$s3 = $p.Slides.AddSlide(3, $layout)
